$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
# Apply header style (bold, centered, bordered) to the new columns K1:O1
# by copying formatting from an existing header cell (A1) before writing values
# (copy happens first so it does not clobber the values we set after).
$ws.Range("A1").Copy() | Out-Null
$ws.Range("K1:O1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A1").Value = "#_Agents"
$ws.Range("B1").Value = "Coverage"
$ws.Range("C1").Value = "Avg_Total_Rounds"
$ws.Range("D1").Value = "Avg_Expl_Cost"
$ws.Range("E1").Value = "Avg_Expl_Eff"
$ws.Range("F1").Value = "Avg_Round_Time"
$ws.Range("G1").Value = "Avg_Agent_Step_Time"
$ws.Range("H1").Value = "Avg_Experiment_Time"
$ws.Range("I1").Value = "Std_Total_Rounds"
$ws.Range("J1").Value = "Std_Expl_Cost"
$ws.Range("K1").Value = "Std_Expl_Eff"
$ws.Range("L1").Value = "Std_Round_Time"
$ws.Range("M1").Value = "Std_Agent_Step_Time"
$ws.Range("N1").Value = "Std_Experiment_Time"
$ws.Range("O1").Value = "Obs_Prob"

# --- Update data rows (rows 2-13) ---
# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 55.488
$ws.Range("D2").Value = 55.488
$ws.Range("E2").Value = 3.081831460000001
$ws.Range("F2").Value = 0.19728432
$ws.Range("G2").Value = 0.19728432
$ws.Range("H2").Value = 10.7967607
$ws.Range("I2").Value = 6.106753710586626
$ws.Range("J2").Value = 6.106753710586626
$ws.Range("K2").Value = 0.3319695128221191
$ws.Range("L2").Value = 0.0271525500770918
$ws.Range("M2").Value = 0.0271525500770918
$ws.Range("N2").Value = 0.6362806974923195
$ws.Range("O2").Value = 0.15

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 91.286
$ws.Range("D3").Value = 91.286
$ws.Range("E3").Value = 1.88692834
$ws.Range("F3").Value = 0.12836258
$ws.Range("G3").Value = 0.12836258
$ws.Range("H3").Value = 11.52921134
$ws.Range("I3").Value = 12.74806771479199
$ws.Range("J3").Value = 12.74806771479199
$ws.Range("K3").Value = 0.2615647180077185
$ws.Range("L3").Value = 0.01844984568920594
$ws.Range("M3").Value = 0.01844984568920594
$ws.Range("N3").Value = 0.9719034653521571
$ws.Range("O3").Value = 0.85

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 31.082
$ws.Range("D4").Value = 62.142
$ws.Range("E4").Value = 2.81263116
$ws.Range("F4").Value = 0.25069254
$ws.Range("G4").Value = 0.12534656
$ws.Range("H4").Value = 3.78075856
$ws.Range("I4").Value = 5.726855135192348
$ws.Range("J4").Value = 11.44919228627138
$ws.Range("K4").Value = 0.5208094055657078
$ws.Range("L4").Value = 0.05354433268679467
$ws.Range("M4").Value = 0.02677210871366603
$ws.Range("N4").Value = 0.4890112879869787
$ws.Range("O4").Value = 0.15

# Row 5
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 48.958
$ws.Range("D5").Value = 95.748
$ws.Range("E5").Value = 1.80630206
$ws.Range("F5").Value = 0.18188218
$ws.Range("G5").Value = 0.09094124
$ws.Range("H5").Value = 4.36648116
$ws.Range("I5").Value = 8.290579196100452
$ws.Range("J5").Value = 14.58224192499625
$ws.Range("K5").Value = 0.2780525084527564
$ws.Range("L5").Value = 0.03029797984375898
$ws.Range("M5").Value = 0.01514911378792565
$ws.Range("N5").Value = 0.5599288927137496
$ws.Range("O5").Value = 0.85

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 16.072
$ws.Range("D6").Value = 64.194
$ws.Range("E6").Value = 2.82603974
$ws.Range("F6").Value = 0.30619246
$ws.Range("G6").Value = 0.07654820000000001
$ws.Range("H6").Value = 1.18651302
$ws.Range("I6").Value = 4.373221482187305
$ws.Range("J6").Value = 17.47177408792705
$ws.Range("K6").Value = 0.7497863354098065
$ws.Range("L6").Value = 0.07938490890889369
$ws.Range("M6").Value = 0.0198462128247929
$ws.Range("N6").Value = 0.307956078591351
$ws.Range("O6").Value = 0.15

# Row 7
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 25.668
$ws.Range("D7").Value = 94.61
$ws.Range("E7").Value = 1.84665818
$ws.Range("F7").Value = 0.24875764
$ws.Range("G7").Value = 0.06218938
$ws.Range("H7").Value = 1.55814086
$ws.Range("I7").Value = 5.658074195979974
$ws.Range("J7").Value = 17.45758960652184
$ws.Range("K7").Value = 0.3382187970593238
$ws.Range("L7").Value = 0.04917224796539548
$ws.Range("M7").Value = 0.01229308110693994
$ws.Range("N7").Value = 0.3059803606889093
$ws.Range("O7").Value = 0.85

# Row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 9.568
$ws.Range("D8").Value = 57.328
$ws.Range("E8").Value = 3.12901554
$ws.Range("F8").Value = 0.36039378
$ws.Range("G8").Value = 0.06006552
$ws.Range("H8").Value = 0.56608716
$ws.Range("I8").Value = 2.471759533284086
$ws.Range("J8").Value = 14.80617647211619
$ws.Range("K8").Value = 0.7522617862617028
$ws.Range("L8").Value = 0.09281849256390275
$ws.Range("M8").Value = 0.01546983806343874
$ws.Range("N8").Value = 0.1828744873046765
$ws.Range("O8").Value = 0.15

# Row 9
$ws.Range("A9").Value = 6
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 17.046
$ws.Range("D9").Value = 87.492
$ws.Range("E9").Value = 1.9961643
$ws.Range("F9").Value = 0.27764202
$ws.Range("G9").Value = 0.04627342
$ws.Range("H9").Value = 0.7736246600000001
$ws.Range("I9").Value = 4.311968244349738
$ws.Range("J9").Value = 15.76283190735342
$ws.Range("K9").Value = 0.369739979340015
$ws.Range("L9").Value = 0.05677194239592664
$ws.Range("M9").Value = 0.009461890468933372
$ws.Range("N9").Value = 0.2043372806911281
$ws.Range("O9").Value = 0.85

# Row 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 7.492
$ws.Range("D10").Value = 59.832
$ws.Range("E10").Value = 3.0224106
$ws.Range("F10").Value = 0.3567655
$ws.Range("G10").Value = 0.04459582
$ws.Range("H10").Value = 0.32952572
$ws.Range("I10").Value = 2.054843245000223
$ws.Range("J10").Value = 16.37645115110747
$ws.Range("K10").Value = 0.7743537773981486
$ws.Range("L10").Value = 0.09786249904938245
$ws.Range("M10").Value = 0.01223290309631887
$ws.Range("N10").Value = 0.1159622019947212
$ws.Range("O10").Value = 0.15

# Row 11
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 13.184
$ws.Range("D11").Value = 82.746
$ws.Range("E11").Value = 2.11282
$ws.Range("F11").Value = 0.2747814
$ws.Range("G11").Value = 0.03434772000000001
$ws.Range("H11").Value = 0.4470119
$ws.Range("I11").Value = 3.526899778004498
$ws.Range("J11").Value = 15.46064510601638
$ws.Range("K11").Value = 0.3934171375721588
$ws.Range("L11").Value = 0.05870578507816147
$ws.Range("M11").Value = 0.007338295412472943
$ws.Range("N11").Value = 0.1382172649985426
$ws.Range("O11").Value = 0.85

# Row 12
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 5.806
$ws.Range("D12").Value = 57.842
$ws.Range("E12").Value = 3.13240196
$ws.Range("F12").Value = 0.34044762
$ws.Range("G12").Value = 0.03404492
$ws.Range("H12").Value = 0.19987644
$ws.Range("I12").Value = 1.54049025327682
$ws.Range("J12").Value = 15.3477079813424
$ws.Range("K12").Value = 0.8463506140142173
$ws.Range("L12").Value = 0.1040565023701384
$ws.Range("M12").Value = 0.01040553304879953
$ws.Range("N12").Value = 0.09040489689664043
$ws.Range("O12").Value = 0.15

# Row 13
$ws.Range("A13").Value = 10
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 10.828
$ws.Range("D13").Value = 75.612
$ws.Range("E13").Value = 2.31506942
$ws.Range("F13").Value = 0.2638943
$ws.Range("G13").Value = 0.02638946
$ws.Range("H13").Value = 0.28237776
$ws.Range("I13").Value = 3.371817656489928
$ws.Range("J13").Value = 14.17431756936165
$ws.Range("K13").Value = 0.4425117447521034
$ws.Range("L13").Value = 0.05393090090342001
$ws.Range("M13").Value = 0.005392926814602189
$ws.Range("N13").Value = 0.09603905325482215
$ws.Range("O13").Value = 0.85

